$wb = $excel.ActiveWorkbook

# The "Status" column on each sheet shrinks because the new status text
# ("Ready for handoff") is shorter than the old text
# ("Handed back: in sync with en-US"), so Excel's column AutoFit ends up
# with a narrower best-fit column width after the values below are written.
$statusColumnWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-17 06:52:58"
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 06:52:52"
$wsZhCn.Columns.Item(3).AutoFit()
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 06:52:58"
$wsDeDe.Columns.Item(3).AutoFit()
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
